$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Apply the updated coin values. For Price (column D) entries that look
# like plain numbers, force Text storage first (NumberFormat "@") so Excel
# doesn't silently convert them to the Number type -- the source data keeps
# these as text (e.g. to preserve trailing zeros).
$ws.Range('D2').Value = '42.296.30'
$ws.Range('E2').Value = '  +0.46%  '
$ws.Range('D3').Value = '2.301.44'
$ws.Range('E3').Value = '  +0.10%  '
$ws.Range('E4').Value = '  -0.07%  '
$ws.Range('D5').NumberFormat = "@"
$ws.Range('D5').Value = '316.58'
$ws.Range('E5').Value = '  +1.33%  '
$ws.Range('D6').NumberFormat = "@"
$ws.Range('D6').Value = '103.28'
$ws.Range('E6').Value = '  -1.06%  '
$ws.Range('D7').NumberFormat = "@"
$ws.Range('D7').Value = '0.630'
$ws.Range('E7').Value = '  +0.65%  '
$ws.Range('E8').Value = '  -0.03%  '
$ws.Range('D9').NumberFormat = "@"
$ws.Range('D9').Value = '0.606'
$ws.Range('E9').Value = '  +0.12%  '
$ws.Range('D10').NumberFormat = "@"
$ws.Range('D10').Value = '39.60'
$ws.Range('E10').Value = '  -1.54%  '
$ws.Range('D11').NumberFormat = "@"
$ws.Range('D11').Value = '0.0909'
$ws.Range('E11').Value = '  -0.33%  '
$ws.Range('D12').NumberFormat = "@"
$ws.Range('D12').Value = '8.37'
$ws.Range('E12').Value = '  +1.31%  '
$ws.Range('E13').Value = '  +0.51%  '
$ws.Range('D14').NumberFormat = "@"
$ws.Range('D14').Value = '0.962'
$ws.Range('E14').Value = '  -1.01%  '
$ws.Range('D15').NumberFormat = "@"
$ws.Range('D15').Value = '15.24'
$ws.Range('E15').Value = '  -2.16%  '
$ws.Range('D16').Value = '2.650.50'
$ws.Range('E16').Value = '  +0.13%  '
$ws.Range('D17').Value = '2.314.21'
$ws.Range('E17').Value = '  +0.73%  '
$ws.Range('D18').Value = '42.416.32'
$ws.Range('E18').Value = '  +0.86%  '
$ws.Range('D19').NumberFormat = "@"
$ws.Range('D19').Value = '7.45'
$ws.Range('E19').Value = '  -2.08%  '
$ws.Range('E20').Value = '  +0.90%  '
$ws.Range('D21').NumberFormat = "@"
$ws.Range('D21').Value = '73.46'
$ws.Range('E21').Value = '  -1.36%  '
$ws.Range('B22').Value = 'PancakeSwap'
$ws.Range('C22').Value = 'https://coinranking.com/coin/ncYFcP709+pancakeswap-cake'
$ws.Range('D22').NumberFormat = "@"
$ws.Range('D22').Value = '3.55'
$ws.Range('E22').Value = '  +3.07%  '
$ws.Range('B23').Value = 'BitcoinCash'
$ws.Range('C23').Value = 'https://coinranking.com/coin/ZlZpzOJo43mIo+bitcoincash-bch'
$ws.Range('D23').NumberFormat = "@"
$ws.Range('D23').Value = '276.77'
$ws.Range('E23').Value = '  +7.52%  '
$ws.Range('D24').NumberFormat = "@"
$ws.Range('D24').Value = '11.38'
$ws.Range('E24').Value = '  +22.30%  '
$ws.Range('D25').NumberFormat = "@"
$ws.Range('D25').Value = '2.27'
$ws.Range('E25').Value = '  -0.85%  '
$ws.Range('E26').Value = '  -0.25%  '
$ws.Range('D27').NumberFormat = "@"
$ws.Range('D27').Value = '10.83'
$ws.Range('E27').Value = '  -1.20%  '
$ws.Range('D29').NumberFormat = "@"
$ws.Range('D29').Value = '22.74'
$ws.Range('E29').Value = '  -0.15%  '
$ws.Range('D30').NumberFormat = "@"
$ws.Range('D30').Value = '37.06'
$ws.Range('E30').Value = '  +3.83%  '
$ws.Range('D31').NumberFormat = "@"
$ws.Range('D31').Value = '165.83'
$ws.Range('E31').Value = '  -0.54%  '
$ws.Range('D32').NumberFormat = "@"
$ws.Range('D32').Value = '0.0875'
$ws.Range('E32').Value = '  -2.17%  '
$ws.Range('E33').Value = '  +1.08%  '
$ws.Range('E34').Value = '  +4.64%  '
$ws.Range('E35').Value = '  -0.41%  '
$ws.Range('D36').NumberFormat = "@"
$ws.Range('D36').Value = '2.61'
$ws.Range('E36').Value = '  -10.47%  '
$ws.Range('D37').NumberFormat = "@"
$ws.Range('D37').Value = '0.0366'
$ws.Range('E37').Value = '  +3.73%  '
$ws.Range('E38').Value = '  +0.80%  '
$ws.Range('D39').NumberFormat = "@"
$ws.Range('D39').Value = '3.70'
$ws.Range('E39').Value = '  +2.00%  '
$ws.Range('E40').Value = '  +0.01%  '
$ws.Range('E41').Value = '  +3.26%  '
$ws.Range('D42').NumberFormat = "@"
$ws.Range('D42').Value = '69.98'
$ws.Range('E42').Value = '  -2.61%  '
$ws.Range('B43').Value = 'Algorand'
$ws.Range('C43').Value = 'https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo'
$ws.Range('D43').NumberFormat = "@"
$ws.Range('D43').Value = '0.227'
$ws.Range('E43').Value = '  -0.20%  '
$ws.Range('B44').Value = 'BitcoinSV'
$ws.Range('C44').Value = 'https://coinranking.com/coin/VcMY11NONHSA0+bitcoinsv-bsv'
$ws.Range('D44').NumberFormat = "@"
$ws.Range('D44').Value = '94.55'
$ws.Range('E44').Value = '  -4.19%  '
$ws.Range('E45').Value = '  +0.04%  '
$ws.Range('D46').NumberFormat = "@"
$ws.Range('D46').Value = '81.07'
$ws.Range('E46').Value = '  +8.85%  '
$ws.Range('D47').NumberFormat = "@"
$ws.Range('D47').Value = '12.05'
$ws.Range('E47').Value = '  -2.12%  '
$ws.Range('D48').NumberFormat = "@"
$ws.Range('D48').Value = '113.00'
$ws.Range('E48').Value = '  +0.82%  '
$ws.Range('E49').Value = '  -0.79%  '
$ws.Range('D50').NumberFormat = "@"
$ws.Range('D50').Value = '5.25'
$ws.Range('E50').Value = '  -1.41%  '
$ws.Range('D51').Value = '1.588.39'
$ws.Range('E51').Value = '  +1.18%  '
